# "added combined complaint/charging step"
#
# Shifts a bunch of shapes/connectors horizontally (to make room) and turns
# the single-line "DA files" / "DA declines" labels into two-line labels
# that call out the new combined complaint/charging step, resizing those
# two text boxes (spAutoFit) to fit.
#
# NOTE on unit conversion: Shape.Left/Top/Width/Height are in points while
# the OOXML stores EMU (1 pt = 12700 EMU). The host's float math is
# single-precision, so a plain "$emu/12700" can come back a hair low and
# truncate to one EMU short when re-emitted. Nudging by +0.5 EMU before
# converting keeps every value landing on the intended integer EMU.

function EmuToPt($emu) {
    return ($emu + 0.5) / 12700
}

function Set-ShapeGeometry($shape, $x, $y, $cx, $cy) {
    if ($null -ne $x) { $shape.Left = EmuToPt $x }
    if ($null -ne $y) { $shape.Top = EmuToPt $y }
    if ($null -ne $cx) { $shape.Width = EmuToPt $cx }
    if ($null -ne $cy) { $shape.Height = EmuToPt $cy }
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Left-hand ("Appear in MACR" / booked-released-referred) cluster ---

# Rectangle 3 "Appear in MACR"
Set-ShapeGeometry ($s.Shapes.Item(1)) 244553 2684725 1796902 1137683

# Rectangle 4 "ACHS"
Set-ShapeGeometry ($s.Shapes.Item(2)) 3432548 843520 1796902 1137683

# Rectangle 5 "JCPSS, state agency, Federal agency"
Set-ShapeGeometry ($s.Shapes.Item(3)) 3572923 4524635 1796902 1598428

# Oval 6 "data ends"
Set-ShapeGeometry ($s.Shapes.Item(4)) 3838357 2464098 1616149 1578935

# Straight Arrow Connector 8 (MACR -> ACHS)
Set-ShapeGeometry ($s.Shapes.Item(5)) 2041455 1412362 1391093 1841205

# Straight Arrow Connector 9 (MACR -> "data ends" oval)
Set-ShapeGeometry ($s.Shapes.Item(6)) 2041455 3253566 1796902 1

# Straight Arrow Connector 11 (MACR -> JCPSS rectangle)
Set-ShapeGeometry ($s.Shapes.Item(7)) 2041455 3253567 1531468 2070282

# TextBox 39 "booked"
Set-ShapeGeometry ($s.Shapes.Item(16)) 2167165 2093272 1139671 461665

# TextBox 40 "released"
Set-ShapeGeometry ($s.Shapes.Item(17)) 2315842 3018768 1272849 461665

# TextBox 42 "referred"
Set-ShapeGeometry ($s.Shapes.Item(18)) 2167165 4024441 1226490 461665

# --- Middle ("ACHS" -> "court records") cluster ---

# Rectangle 13 "court records"
Set-ShapeGeometry ($s.Shapes.Item(8)) 7344527 843520 1796902 1137683

# Straight Arrow Connector 15 (ACHS -> court records): also loses its
# vertical flip now that it is perfectly horizontal (cy = 0).
$connector15 = $s.Shapes.Item(9)
Set-ShapeGeometry $connector15 5229450 1412362 2115077 0
$connector15.VerticalFlip = 0

# Oval 17 "data ends" (ACHS branch)
Set-ShapeGeometry ($s.Shapes.Item(10)) 7434903 2481817 1616149 1578935

# Straight Arrow Connector 18 (ACHS -> "data ends" oval)
Set-ShapeGeometry ($s.Shapes.Item(11)) 5229450 1412362 2205453 1858923

# TextBox 43: "DA files" -> two-line "complaint " / "& charged"
$daFiles = $s.Shapes.Item(19)
Set-ShapeGeometry $daFiles 5648802 887939 1482768 830997
$daFiles.TextFrame.WordWrap = -1
$daFiles.TextFrame.TextRange.Text = "complaint `r& charged"
$daFiles.TextFrame.TextRange.ParagraphFormat.Alignment = 2

# TextBox 44: "DA declines" -> two-line "complaint &" / "DA declines"
$daDeclines = $s.Shapes.Item(20)
Set-ShapeGeometry $daDeclines 5465022 2025509 1786386 830997
$daDeclines.TextFrame.TextRange.Text = "complaint &`rDA declines"
$daDeclines.TextFrame.TextRange.ParagraphFormat.Alignment = 2

# --- Right ("court records" -> open/sealed) cluster ---

# Straight Arrow Connector 22 (court records -> "open" branch)
Set-ShapeGeometry ($s.Shapes.Item(12)) 9141429 1412362 1449571 0

# Oval 23 "not public"
Set-ShapeGeometry ($s.Shapes.Item(13)) 10345479 2441008 1616149 1578935

# Straight Arrow Connector 24 (court records -> "not public" oval)
Set-ShapeGeometry ($s.Shapes.Item(14)) 9141429 1412362 1204050 1818114

# Oval 27 "public"
Set-ShapeGeometry ($s.Shapes.Item(15)) 10345480 599801 1616149 1578935

# TextBox 45 "open"
Set-ShapeGeometry ($s.Shapes.Item(21)) 9345175 1140718 835485 461665

# TextBox 46 "sealed"
Set-ShapeGeometry ($s.Shapes.Item(22)) 9360306 2185597 1011815 461665
